$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Start Date (A2) to the new date value (2019-05-24 -> serial 43609)
$ws.Range("A2").Value = (Get-Date -Year 2019 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0).Date

# Clear End Date (B2) - leave it blank
$ws.Range("B2").ClearContents()

# Update the active selection to F3
$ws.Range("F3").Select()
